{"js": "// Office.js (Word JavaScript API) script.\n//\n// The document currently ends with one paragraph:\n//   \"Edit: Fout gevonden! ... neerzetten. \" + [_GoBack bookmark]\n//\n// Target: keep that paragraph's text as-is, drop the bookmark from it,\n// append the \"18 December\" / \"21 December\" log entries as new\n// paragraphs, and finish with a brand-new (empty) last paragraph that\n// now owns the relocated \"_GoBack\" bookmark.\n\nconst doc = context.document;\nconst body = doc.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// --- 1. Detach the _GoBack bookmark from the last paragraph -----------\ndoc.deleteBookmark(\"_GoBack\");\n\n// --- 2. Insert all the new paragraphs (no bold yet, to avoid bold\n//        \"bleeding\" into paragraphs inserted after a bold one) ---------\nconst blank1 = lastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nconst heading18 = blank1.insertParagraph(\"18 December:\", Word.InsertLocation.after);\nconst collisionPara = heading18.insertParagraph(\n  \"Ik heb de eerste collision checks gemaakt. Via een speciale class die methodes heeft om te controleren of objecten botsen. De spelen kan nu niet meer door de vloer vallen, behalve op plekken waar geen vloer zit. \",\n  Word.InsertLocation.after\n);\nconst objectsPara = collisionPara.insertParagraph(\n  \"Ik heb ook al een test gedaan met opbjecten, maar deze zakken toch nog door de vloer heen.\",\n  Word.InsertLocation.after\n);\nconst blank2 = objectsPara.insertParagraph(\"\", Word.InsertLocation.after);\nconst heading21 = blank2.insertParagraph(\"21 December:\", Word.InsertLocation.after);\nconst finalParagraph = heading21.insertParagraph(\"\", Word.InsertLocation.after);\n\n// --- 3. Bold only the two date-heading paragraphs -----------------------\nheading18.font.bold = true;\nheading21.font.bold = true;\n\n// --- 4. Re-create the _GoBack bookmark in the new, final, empty paragraph\nconst bookmarkRange = finalParagraph.getRange();\nbookmarkRange.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Word COM interop script\n#\n# Source document ends with one paragraph:\n#   \"Edit: Fout gevonden! ... neerzetten. \" + [_GoBack bookmark]\n#\n# Target: keep that paragraph's text as-is, drop the bookmark from it,\n# append the \"18 December\" / \"21 December\" log entries as new paragraphs,\n# and finish with a brand-new (empty) last paragraph that now owns the\n# relocated \"_GoBack\" bookmark.\n\n$d = $word.ActiveDocument\n\n# --- 1. Detach the _GoBack bookmark from the last paragraph -----------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- 2. Append all the new paragraphs in one shot ----------------------\n# Each backtick-r (`r) is a paragraph mark. Text before the first one is\n# appended to the existing last paragraph (it is empty here, so that\n# paragraph's text is left untouched).\n$last = $d.Paragraphs.Last\n$insertionRange = $last.Range\n$insertionRange.Collapse(0)   # wdCollapseEnd\n\n$newText = \"`r`r18 December:`rIk heb de eerste collision checks gemaakt. Via een speciale class die methodes heeft om te controleren of objecten botsen. De spelen kan nu niet meer door de vloer vallen, behalve op plekken waar geen vloer zit. `rIk heb ook al een test gedaan met opbjecten, maar deze zakken toch nog door de vloer heen.`r`r21 December:`r\"\n\n$insertionRange.InsertAfter($newText)\n\n# --- 3. Bold the two date-heading paragraphs ---------------------------\n# Use Paragraph.Range (which, unlike the Office.js Range, includes the\n# trailing paragraph mark) so both the run AND the paragraph mark pick up\n# the Bold run property, matching how Word itself records a fully-bold\n# paragraph.\n$paragraphs = $d.Paragraphs\n$count = $paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -eq \"18 December:`r\" -or $t -eq \"21 December:`r\") {\n        $p.Range.Font.Bold = 1\n    }\n}\n\n# --- 4. Re-create the _GoBack bookmark in the new, final, empty paragraph\n$finalParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$bookmarkRange = $finalParagraph.Range\n$bookmarkRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n"}
